# Saldo.xlsx ("Export" sheet) update:
#   - Add a new account row (004451978 / ANTONIO / 4203.98) right above the
#     004392159 / RODRIGO row (i.e. as the new row 4, pushing everything
#     below it down by one row).
#   - Remove the old 004451978 / ANTONIO / 12.17 row that used to sit
#     between the 004451652/MATEUS (12.26) and 004809902/PEDRO (12.16) rows.
#     After the insertion above, that row has shifted down from row 152 to
#     row 153.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row ---------------------------------------------------
$ws.Rows.Item(4).Insert()

# Column A holds zero-padded account numbers that must stay TEXT (not be
# coerced to a number, which would drop the leading zero). Assigning a
# numeric-looking string straight to .Value/.Value2/.Formula gets
# auto-converted to a number by this host, and forcing a "@" text
# NumberFormat leaves a stray style on the cell. Instead, compute the text
# with a throwaway formula cell and paste-special just the resulting value,
# which preserves both the string type and the plain/default cell style.
$helper = $ws.Cells.Item(1, 10)
$helper.Formula = '="004451978"'
$helper.Copy()
$ws.Cells.Item(4, 1).PasteSpecial(-4163)  # xlPasteValues
$helper.ClearContents()

$ws.Cells.Item(4, 2).Value = "ANTONIO"
$ws.Cells.Item(4, 3).Value = 4203.98

# --- Remove the old, now-duplicate row ------------------------------------
$ws.Rows.Item(153).Delete()
